$wb = $excel.ActiveWorkbook

# --- Sheet: status --- add row 3 (duplicate of row 2) ---
$ws = $wb.Worksheets.Item("status")
$ws.Cells.Item(3,1).Value = 1465993372
$ws.Cells.Item(3,2).Value = 1041322197

# --- Sheet: neighbors --- add row 3 (duplicate of row 2) ---
$ws = $wb.Worksheets.Item("neighbors")
$ws.Cells.Item(3,1).Value = 1465993372
$ws.Cells.Item(3,2).Value = "10.0.0.5"
$ws.Cells.Item(3,3).Value = $true
$ws.Cells.Item(3,4).Value = $true
$ws.Cells.Item(3,5).Value = $false
$ws.Cells.Item(3,6).Value = 3
$ws.Cells.Item(3,7).Value = 1

# --- Sheet: links --- add row 3 (duplicate of row 2) ---
$ws = $wb.Worksheets.Item("links")
$ws.Cells.Item(3,1).Value = 1465993372
$ws.Cells.Item(3,2).Value = "10.0.0.6"
$ws.Cells.Item(3,3).Value = "10.0.0.5"
$ws.Cells.Item(3,4).Value = 37752
$ws.Cells.Item(3,5).Value = 0.886
$ws.Cells.Item(3,6).Value = 0.886
$ws.Cells.Item(3,7).Value = 1303

# --- Sheet: routes --- add rows 5,6,7 (duplicate of rows 2,3,4) ---
$ws = $wb.Worksheets.Item("routes")
$ws.Cells.Item(5,1).Value = 1465993372
$ws.Cells.Item(5,2).Value = "10.0.0.3"
$ws.Cells.Item(5,3).Value = 32
$ws.Cells.Item(5,4).Value = "10.0.0.5"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 4767
$ws.Cells.Item(5,7).Value = "mesh0"

$ws.Cells.Item(6,1).Value = 1465993372
$ws.Cells.Item(6,2).Value = "10.0.0.4"
$ws.Cells.Item(6,3).Value = 32
$ws.Cells.Item(6,4).Value = "10.0.0.5"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 2327
$ws.Cells.Item(6,7).Value = "mesh0"

$ws.Cells.Item(7,1).Value = 1465993372
$ws.Cells.Item(7,2).Value = "10.0.0.5"
$ws.Cells.Item(7,3).Value = 32
$ws.Cells.Item(7,4).Value = "10.0.0.5"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 1303
$ws.Cells.Item(7,7).Value = "mesh0"

# --- Sheet: topology --- add rows 8-13 (duplicate of rows 2-7) ---
$ws = $wb.Worksheets.Item("topology")
$ws.Cells.Item(8,1).Value = 1465993372
$ws.Cells.Item(8,2).Value = "10.0.0.4"
$ws.Cells.Item(8,3).Value = "10.0.0.3"
$ws.Cells.Item(8,4).Value = 1
$ws.Cells.Item(8,5).Value = 0.274
$ws.Cells.Item(8,6).Value = 3730
$ws.Cells.Item(8,7).Value = 262741

$ws.Cells.Item(9,1).Value = 1465993372
$ws.Cells.Item(9,2).Value = "10.0.0.3"
$ws.Cells.Item(9,3).Value = "10.0.0.4"
$ws.Cells.Item(9,4).Value = 0.419
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 2440
$ws.Cells.Item(9,7).Value = 281330

$ws.Cells.Item(10,1).Value = 1465993372
$ws.Cells.Item(10,2).Value = "10.0.0.5"
$ws.Cells.Item(10,3).Value = "10.0.0.4"
$ws.Cells.Item(10,4).Value = 1
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 1024
$ws.Cells.Item(10,7).Value = 281330

$ws.Cells.Item(11,1).Value = 1465993372
$ws.Cells.Item(11,2).Value = "10.0.0.4"
$ws.Cells.Item(11,3).Value = "10.0.0.5"
$ws.Cells.Item(11,4).Value = 1
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 1024
$ws.Cells.Item(11,7).Value = 280263

$ws.Cells.Item(12,1).Value = 1465993372
$ws.Cells.Item(12,2).Value = "10.0.0.6"
$ws.Cells.Item(12,3).Value = "10.0.0.5"
$ws.Cells.Item(12,4).Value = 0.886
$ws.Cells.Item(12,5).Value = 0.886
$ws.Cells.Item(12,6).Value = 1303
$ws.Cells.Item(12,7).Value = 280263

$ws.Cells.Item(13,1).Value = 1465993372
$ws.Cells.Item(13,2).Value = "10.0.0.5"
$ws.Cells.Item(13,3).Value = "10.0.0.6"
$ws.Cells.Item(13,4).Value = 0.886
$ws.Cells.Item(13,5).Value = 0.886
$ws.Cells.Item(13,6).Value = 1303
$ws.Cells.Item(13,7).Value = 0
